# Setting up run modes for Test data (Skipped Test Case for Add Customer Test)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddCustomerTest")

# New "Runmode" column (E) next to the existing data (A:D)
$ws.Range("E1").Value = "Runmode"
$ws.Range("E2").Value = "Y"
$ws.Range("E3").Value = "N"   # Skipped Test Case for Add Customer Test
$ws.Range("E4").Value = "Y"
$ws.Range("E5").Value = "Y"
$ws.Range("E6").Value = "Y"
$ws.Range("E7").Value = "Y"

# Move the selection to the last populated cell on this sheet
$ws.Range("E7").Select()

# Make the "AddCustomerTest" sheet the active/selected tab
$ws.Activate()
